$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new price-history row (row 82) exactly like the existing
# trailing rows: both the date and the value are stored as plain text,
# not as a date-serial number / float.
$ws.Cells.Item(82, 1).Value = "'2025-01-19"
$ws.Cells.Item(82, 2).Value = "'42.6"

# Drop the "quote prefix" style Excel auto-applies for the leading
# apostrophe so the new cells fall back to the default (unstyled) format,
# matching the rest of the sheet's plain text rows.
$ws.Cells.Item(82, 1).Style = "Normal"
$ws.Cells.Item(82, 2).Style = "Normal"
